$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 14 (existing rows 14-29 shift down to 17-32),
# carrying formatting down the way Excel normally does.
$ws.Rows("14:16").Insert()

# Row 14: new Castle Brite / Primera record ($/caja 15 kilos)
$ws.Cells.Item(14, 1).Value = 9
$ws.Cells.Item(14, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44533
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100103
$ws.Cells.Item(14, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(14, 9).Value = 100103003
$ws.Cells.Item(14, 10).Value = "Damasco"
$ws.Cells.Item(14, 11).Value = "Castle Brite"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 410
$ws.Cells.Item(14, 14).Value = 16500
$ws.Cells.Item(14, 15).Value = 16500
$ws.Cells.Item(14, 16).Value = 16500
$ws.Cells.Item(14, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(14, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(14, 19).Value = 1100
$ws.Cells.Item(14, 20).Value = 15

# Row 15: new Castle Brite / Segunda record ($/caja 15 kilos)
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44533
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103003
$ws.Cells.Item(15, 10).Value = "Damasco"
$ws.Cells.Item(15, 11).Value = "Castle Brite"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 300
$ws.Cells.Item(15, 14).Value = 10500
$ws.Cells.Item(15, 15).Value = 10500
$ws.Cells.Item(15, 16).Value = 10500
$ws.Cells.Item(15, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(15, 19).Value = 700
$ws.Cells.Item(15, 20).Value = 15

# Row 16: new Castle Brite / Tercera record ($/caja 15 kilos)
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44533
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100103
$ws.Cells.Item(16, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(16, 9).Value = 100103003
$ws.Cells.Item(16, 10).Value = "Damasco"
$ws.Cells.Item(16, 11).Value = "Castle Brite"
$ws.Cells.Item(16, 12).Value = "Tercera"
$ws.Cells.Item(16, 13).Value = 250
$ws.Cells.Item(16, 14).Value = 7500
$ws.Cells.Item(16, 15).Value = 7500
$ws.Cells.Item(16, 16).Value = 7500
$ws.Cells.Item(16, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(16, 19).Value = 500
$ws.Cells.Item(16, 20).Value = 15
